$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Content edit: falafel options text changes from "Lentil,Broccoli,Split pea"
# to "Broccoli,Split pea" on Hoja2 (this also re-packs the shared-string table,
# which is what shifts every other shared-string index referenced elsewhere). ---
$ws2.Range("C10").Value = "Broccoli,Split pea"

# --- Hoja1 rows 4-13: populate with the same trivia rows that live on Hoja2
# (rows 1-10), but dated 30-08-2023 (serial 45168) instead of 26-08-2023. ---

# Row 4 <- Hoja2 row 1
$ws1.Range("A4").Value = "What 90s boy band member bought Myspace in 2011?"
$ws1.Range("B4").Value = "Justin Timberlake"
$ws1.Range("C4").Value = "Nick Lachey,Shawn Stockman,AJ McLean"
$ws1.Range("D4").Value = 45168

# Row 5 <- Hoja2 row 2
$ws1.Range("A5").Value = "What's the name of Hagrid's pet spider?"
$ws1.Range("B5").Value = "Aragog"
$ws1.Range("C5").Value = "Nigini,Crookshanks,Mosag"
$ws1.Range("D5").Value = 45168

# Row 6 <- Hoja2 row 3
$ws1.Range("A6").Value = "What's the heaviest organ in the human body?"
$ws1.Range("B6").Value = "Liver"
$ws1.Range("C6").Value = "Brain,Skin,Heart"
$ws1.Range("D6").Value = 45168

# Row 7 <- Hoja2 row 4
$ws1.Range("A7").Value = "Which of these EU countries does not use the euro as its currency?"
$ws1.Range("B7").Value = " All are correct"
$ws1.Range("C7").Value = "Poland,Denmark,Sweden"
$ws1.Range("D7").Value = 45168

# Row 8 <- Hoja2 row 5
$ws1.Range("A8").Value = "What element does the chemical symbol Au stand for?"
$ws1.Range("B8").Value = "Gold"
$ws1.Range("C8").Value = "Salt,Magnesium,Silver"
$ws1.Range("D8").Value = 45168

# Row 9 <- Hoja2 row 6
$ws1.Range("A9").Value = "On average, how many seeds are located on the outside of a strawberry?"
$ws1.Range("B9").Value = 200
$ws1.Range("C9").Value = "100,400,500"
$ws1.Range("D9").Value = 45168

# Row 10 <- Hoja2 row 7
$ws1.Range("A10").Value = "What is the oldest soft drink in the United States?"
$ws1.Range("B10").Value = "Dr. Pepper"
$ws1.Range("C10").Value = "Coca Cola,Pepsi,Canada Dry Ginger Ale"
$ws1.Range("D10").Value = 45168

# Row 11 <- Hoja2 row 8 (question has a rich-text run, so copy the cell itself
# to carry the formatted run over rather than retyping plain text).
$ws2.Range("A8").Copy($ws1.Range("A11"))
$ws1.Range("B11").Value = "Central Perk"
$ws1.Range("C11").Value = "Java Park,Central Park Coffee,Central Park Roastery"
$ws1.Range("D11").Value = 45168

# Row 12 <- Hoja2 row 9
$ws1.Range("A12").Value = "Which country's national animal is a unicorn?"
$ws1.Range("B12").Value = "Scotland"
$ws1.Range("C12").Value = "Denmark,New Zealand,France"
$ws1.Range("D12").Value = 45168
# C12 previously used an underlined font (now unused elsewhere) - drop the
# underline so it matches the plain centered style used by the rest of column C.
$ws1.Range("C12").Font.Underline = -4142

# Row 13 <- Hoja2 row 10
$ws1.Range("A13").Value = "What is the main ingredient in a falafel?"
$ws1.Range("B13").Value = "Lentil"
$ws1.Range("C13").Value = "Broccoli,Split pea"
$ws1.Range("D13").Value = 45168

# --- Page setup: Hoja2 gains an explicit print page setup (Letter, portrait). ---
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- Selections: Hoja2 remembers a selection on C10; Hoja1 stays the active
# sheet/tab with its selection moved to D4. ---
$ws2.Activate()
$ws2.Range("C10").Select()

$ws1.Activate()
$ws1.Range("D4").Select()
